$wb = $excel.ActiveWorkbook

# Updated "想去人数" (interest count) values, refreshed from the source site.
$updates = @{
    "F2"  = 224
    "F3"  = 259
    "F4"  = 271
    "F7"  = 6338
    "F8"  = 48
    "F10" = 107
    "F14" = 191
    "F15" = 484
    "F16" = 42
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
